$d = $word.ActiveDocument

# Add a new, centered paragraph with the contact info directly below the
# name ("Dheeraj Chand") paragraph. Using Find/Replace with a paragraph
# mark (^p) in the replacement text creates a brand-new plain paragraph
# that does not inherit the bold/large-font run formatting of the name
# line, while still picking up the centered alignment of that paragraph.
$d.Content.Find.Execute("Dheeraj Chand", $false, $false, $false, $false, `
    $false, $true, 1, $false, `
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX", `
    2)
